$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analysis timestamp (shared string)
$ws.Range("A2").Value = "2025-05-26 17:28:57"

# Update numeric metrics in row 2
$ws.Range("C2").Value = 11756
$ws.Range("D2").Value = 72.43822786370077
$ws.Range("E2").Value = 2236
$ws.Range("F2").Value = 13.77780516359603

$ws.Range("O2").Value = 4695
$ws.Range("P2").Value = 28.92969375808738
$ws.Range("Q2").Value = 497862.3
$ws.Range("R2").Value = 3572
$ws.Range("S2").Value = 22.00998213075359
$ws.Range("T2").Value = 3429
$ws.Range("U2").Value = 21.12884342843059
$ws.Range("V2").Value = 2424263.57
$ws.Range("W2").Value = 2296
$ws.Range("X2").Value = 14.14751371002526

$ws.Range("AI2").Value = 484
$ws.Range("AJ2").Value = 863
$ws.Range("AK2").Value = 1339
$ws.Range("AL2").Value = 18.01935964259122
$ws.Range("AM2").Value = 32.12956068503351
$ws.Range("AN2").Value = 49.85107967237528
$ws.Range("AO2").Value = 1331527.2
$ws.Range("AP2").Value = 249630.26
$ws.Range("AQ2").Value = 83316.29000000001
$ws.Range("AR2").Value = 79.99688790526135
$ws.Range("AS2").Value = 14.99754862460282
$ws.Range("AT2").Value = 5.005563470135831
$ws.Range("AU2").Value = 50.29909706546275
$ws.Range("AV2").Value = 234.8590021691974
$ws.Range("AW2").Value = 552.6344878408254
